$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create the 3 new rows (805-807) by copying the existing row 804
# (template "Especial/Primera/Segunda" Brasil block, D=44454) down, before row 804's
# own values are overwritten in Step 2.
$ws.Range("A804:T804").Copy($ws.Range("A805:T805"))
$ws.Range("A804:T804").Copy($ws.Range("A806:T806"))
$ws.Range("A804:T804").Copy($ws.Range("A807:T807"))

# Row 805: Especial
$ws.Cells.Item(805, 12).Value = "Especial"

# Row 806: Primera
$ws.Cells.Item(806, 12).Value = "Primera"

# Row 807: Segunda (matches the original row-804 content exactly, nothing else to change)

# --- Step 2: apply the updated values to rows 761-804.
# Column map: D=4, L=12, M=13, N=14, O=15, P=16, R=18, S=19
$changes = @(
    @{Row=761; D=44753; N=7000; O=7500; P=7250; R="Brasil"; S=1812}
    @{Row=762; D=44753; N=7000; O=7500; P=7250; R="Brasil"; S=1812}
    @{Row=763; D=44753; N=7000; O=7500; P=7250; R="Brasil"; S=1812}
    @{Row=764; D=44489; N=6500; O=7000; P=6750; S=1688}
    @{Row=765; D=44489; N=6500; O=7000; P=6750; S=1688}
    @{Row=766; D=44489; N=6500; O=7000; P=6750; S=1688}
    @{Row=767; D=44659; N=7000; O=7500; P=7250; S=1812}
    @{Row=768; D=44659; N=7000; O=7500; P=7250; S=1812}
    @{Row=769; D=44659; N=7000; O=7500; P=7250; S=1812}
    @{Row=770; D=44505; N=6500; O=7000; P=6750; S=1688}
    @{Row=771; D=44505; N=6500; O=7000; P=6750; S=1688}
    @{Row=772; D=44505; N=6500; O=7000; P=6750; S=1688}
    @{Row=773; D=44340; M=512; N=8000; O=8500; P=8250; S=2062}
    @{Row=774; D=44340; M=512; N=8000; O=8500; P=8250; S=2062}
    @{Row=775; D=44340; M=512; N=8000; O=8500; P=8250; S=2062}
    @{Row=776; D=44326; M=500; N=9000; O=9500; P=9250; R="Perú"; S=2312}
    @{Row=777; D=44326; M=500; N=9000; O=9500; P=9250; R="Perú"; S=2312}
    @{Row=778; D=44326; M=500; N=9000; O=9500; P=9250; R="Perú"; S=2312}
    @{Row=779; D=44714; M=512; N=10000; O=11000; P=10500; R="Brasil"; S=2625}
    @{Row=780; D=44714; M=512; N=10000; O=11000; P=10500; R="Brasil"; S=2625}
    @{Row=781; D=44714; M=512; N=10000; O=11000; P=10500; R="Brasil"; S=2625}
    @{Row=782; D=44343; L="Especial"; M=400; N=8000; O=8500; P=8250; R="Perú"; S=2062}
    @{Row=783; D=44343; L="Primera"; M=400; N=8000; O=8500; P=8250; R="Perú"; S=2062}
    @{Row=784; D=44343; L="Segunda"; M=512; N=8000; O=8500; P=8250; S=2062}
    @{Row=785; D=44426; N=8500; O=9000; P=8750; R="Brasil"; S=2188}
    @{Row=786; D=44426; N=8500; O=9000; P=8750; R="Brasil"; S=2188}
    @{Row=787; D=44259}
    @{Row=788; D=44259}
    @{Row=789; D=44259}
    @{Row=790; D=44238; N=5000; O=5500; P=5250; S=1312}
    @{Row=791; D=44238; N=5000; O=5500; P=5250; S=1312}
    @{Row=792; D=44238; N=5000; O=5500; P=5250; S=1312}
    @{Row=793; D=44622; N=5500; O=6000; P=5750; R="Perú"; S=1438}
    @{Row=794; D=44622; N=5500; O=6000; P=5750; R="Perú"; S=1438}
    @{Row=795; D=44622; N=5500; O=6000; P=5750; R="Perú"; S=1438}
    @{Row=796; D=44406; N=9000; O=9500; P=9250; R="Brasil"; S=2312}
    @{Row=797; D=44406; N=9000; O=9500; P=9250; R="Brasil"; S=2312}
    @{Row=798; D=44406; N=9000; O=9500; P=9250; R="Brasil"; S=2312}
    @{Row=799; D=44627; N=6500; O=7000; P=6750; S=1688}
    @{Row=800; D=44627; N=6500; O=7000; P=6750; S=1688}
    @{Row=801; D=44627; N=6500; O=7000; P=6750; S=1688}
    @{Row=802; D=44547; N=5500; O=6000; P=5750; R="Perú"; S=1438}
    @{Row=803; D=44547; N=5500; O=6000; P=5750; R="Perú"; S=1438}
    @{Row=804; D=44547; L="Segunda"; M=512; N=5500; O=6000; P=5750; R="Perú"; S=1438}
)

foreach ($item in $changes) {
    $r = $item.Row
    if ($item.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $item.D }
    if ($item.ContainsKey("L")) { $ws.Cells.Item($r, 12).Value = $item.L }
    if ($item.ContainsKey("M")) { $ws.Cells.Item($r, 13).Value = $item.M }
    if ($item.ContainsKey("N")) { $ws.Cells.Item($r, 14).Value = $item.N }
    if ($item.ContainsKey("O")) { $ws.Cells.Item($r, 15).Value = $item.O }
    if ($item.ContainsKey("P")) { $ws.Cells.Item($r, 16).Value = $item.P }
    if ($item.ContainsKey("R")) { $ws.Cells.Item($r, 18).Value = $item.R }
    if ($item.ContainsKey("S")) { $ws.Cells.Item($r, 19).Value = $item.S }
}
